$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-150) holds the "Förändrad" date value.
# The workbook was refreshed a day later, so the value shifts
# from serial 45181 (2023-09-12) to 45182 (2023-09-13).
$oldSerial = 45181
$newSerial = 45182

$lastRow = 150
$range = $ws.Range("C2:C$lastRow")

for ($i = 1; $i -le $range.Rows.Count; $i++) {
    $cell = $range.Cells.Item($i, 1)
    if ($cell.Value2 -eq $oldSerial) {
        $cell.Value = $newSerial
    }
}
